$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add score results for the two newly completed matches (rows 16 and 17)
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 4

$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 3

# Update the active selection to H19 as in the saved workbook state
$ws.Range("H19").Select()
